$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(27).Copy()
$ws.Rows.Item(29).PasteSpecial(-4122)
